# Apply the "clean versions" data fill to Sheet1.
# Rows 25-54 (robustness values 21-50) are missing their B, C, D, E
# (clean-network) values; this script fills them back in. For these rows,
# B/C/D mirror the existing F/G/H (attack) columns, while E holds new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    25 = @(0.99658725935302395, 0.99589069204575498, 0.996568721619426,   0.99650607378946698)
    26 = @(0.99661237757294696, 0.99593539809556397, 0.99664044862235202, 0.99654417917110905)
    27 = @(0.99661991424251195, 0.99588785757863496, 0.99670260684823897, 0.99658331526806898)
    28 = @(0.99667865832259195, 0.99592370683897402, 0.99673761517576498, 0.99662348208034401)
    29 = @(0.996720382940538,   0.99595012572346397, 0.996776119680719,   0.99666260427082698)
    30 = @(0.99674085297646797, 0.99602302510928797, 0.99681020690868305, 0.99669625496928904)
    31 = @(0.99694285317970499, 0.99599643101806801, 0.99683700486621696, 0.99672983011170802)
    32 = @(0.99699267693871196, 0.99601399980604499, 0.99688874209003198, 0.99676880587517303)
    33 = @(0.99698371179771195, 0.996025246055387,   0.99691372957420299, 0.99680409887645305)
    34 = @(0.99700073496099295, 0.996066857505162,   0.99696767681228804, 0.99684484373708004)
    35 = @(0.99704837527734103, 0.99610583081454296, 0.99709948224774203, 0.99688454657873105)
    36 = @(0.99708926440836598, 0.99610908002187404, 0.99714867748664004, 0.99689989360258202)
    37 = @(0.99712520684953398, 0.99619732308507103, 0.997199885223222,   0.99691292724316605)
    38 = @(0.99711753835732597, 0.99622026958930998, 0.99724132243302699, 0.99694020951905704)
    39 = @(0.99715960953885296, 0.99620788709903696, 0.99725481089719503, 0.99696852251026502)
    40 = @(0.99719726991324797, 0.99628777980647998, 0.99735579563933696, 0.99699786621678799)
    41 = @(0.99723274533477801, 0.99633012748148098, 0.99738055162157602, 0.99702824063862805)
    42 = @(0.99724762935476396, 0.99637175365576103, 0.99738428101136001, 0.99705964577578499)
    43 = @(0.99725551759904696, 0.99634234145798195, 0.99741124566927697, 0.99709208162825802)
    44 = @(0.99742890206061996, 0.99642317653371504, 0.99745341010599897, 0.99720301815616996)
    45 = @(0.99744039453639799, 0.99642975593315197, 0.99748263348803601, 0.99722175289723003)
    46 = @(0.99745452024440096, 0.99643180263928,     0.99749035403657305, 0.99728328927880305)
    47 = @(0.99746230459912399, 0.99655619525346895, 0.99749938670206895, 0.997315567028768)
    48 = @(0.99746346292681298, 0.99652310192955895, 0.99751377172495703, 0.997326822276418)
    49 = @(0.99750197895852699, 0.99656554776792405, 0.99754943038476596, 0.99732900760111498)
    50 = @(0.99751057234297102, 0.99656516493080605, 0.99756577785689604, 0.99739968685886904)
    51 = @(0.997552290954417,   0.996689229515756,   0.99760743325708601, 0.99743465636724604)
    52 = @(0.99758592417682601, 0.99669495980209799, 0.99761046895910099, 0.997465352943397)
    53 = @(0.99760436498268201, 0.99669088602251499, 0.99762820217071202, 0.99746621187282802)
    54 = @(0.99768105838250498, 0.99676716795423903, 0.99763178022531096, 0.99748149263640795)
}

foreach ($row in $data.Keys | Sort-Object) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
}

# Update the active selection to match the post-edit state.
$ws.Range("D4:D54").Select()
